# Apply updated probability matrix values to Sheet1.
# These values represent a refreshed team-specific transition matrix
# (added more games, sped up simulate game logic, drafted optimization logic).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

    $ws.Range("B2").Value = 0.165938864628821
    $ws.Range("C2").Value = 0.6026200873362445
    $ws.Range("J2").Value = 0.008733624454148471
    $ws.Range("P2").Value = 0.1266375545851528
    $ws.Range("S2").Value = 0.09606986899563319
    $ws.Range("C3").Value = 0.0072992700729927
    $ws.Range("J3").Value = 0.0072992700729927
    $ws.Range("P3").Value = 0.7883211678832117
    $ws.Range("S3").Value = 0.1970802919708029
    $ws.Range("J4").Value = 0.02777777777777778
    $ws.Range("P4").Value = 0.6944444444444444
    $ws.Range("S4").Value = 0.2777777777777778
    $ws.Range("B6").Value = 0.03448275862068965
    $ws.Range("D6").Value = 0.01970443349753695
    $ws.Range("F6").Value = 0.05911330049261083
    $ws.Range("J6").Value = 0.2167487684729064
    $ws.Range("O6").Value = 0.01970443349753695
    $ws.Range("Q6").Value = 0.1625615763546798
    $ws.Range("R6").Value = 0.08866995073891626
    $ws.Range("S6").Value = 0.3990147783251232
    $ws.Range("B7").Value = 0.1132075471698113
    $ws.Range("D7").Value = 0.03773584905660377
    $ws.Range("F7").Value = 0.06289308176100629
    $ws.Range("J7").Value = 0.119496855345912
    $ws.Range("O7").Value = 0.006289308176100629
    $ws.Range("Q7").Value = 0.1635220125786163
    $ws.Range("R7").Value = 0.07547169811320754
    $ws.Range("S7").Value = 0.4213836477987422
    $ws.Range("B8").Value = 0.06937799043062201
    $ws.Range("D8").Value = 0.01435406698564593
    $ws.Range("F8").Value = 0.05502392344497608
    $ws.Range("J8").Value = 0.1267942583732057
    $ws.Range("O8").Value = 0.01913875598086124
    $ws.Range("Q8").Value = 0.222488038277512
    $ws.Range("R8").Value = 0.09330143540669857
    $ws.Range("S8").Value = 0.3995215311004784
    $ws.Range("B9").Value = 0.1016042780748663
    $ws.Range("D9").Value = 0.0053475935828877
    $ws.Range("F9").Value = 0.1016042780748663
    $ws.Range("J9").Value = 0.0855614973262032
    $ws.Range("O9").Value = 0.0160427807486631
    $ws.Range("Q9").Value = 0.1871657754010695
    $ws.Range("R9").Value = 0.0855614973262032
    $ws.Range("S9").Value = 0.4171122994652406
    $ws.Range("B10").Value = 0.1041666666666667
    $ws.Range("D10").Value = 0.01721014492753623
    $ws.Range("E10").Value = 0.001811594202898551
    $ws.Range("F10").Value = 0.07065217391304347
    $ws.Range("J10").Value = 0.1032608695652174
    $ws.Range("O10").Value = 0.009963768115942028
    $ws.Range("Q10").Value = 0.2327898550724638
    $ws.Range("R10").Value = 0.09963768115942029
    $ws.Range("S10").Value = 0.3605072463768116
    $ws.Range("G11").Value = 0.1486486486486487
    $ws.Range("J11").Value = 0.05855855855855856
    $ws.Range("K11").Value = 0.1801801801801802
    $ws.Range("L11").Value = 0.5990990990990991
    $ws.Range("S11").Value = 0.01351351351351351
    $ws.Range("G12").Value = 0.7338129496402878
    $ws.Range("J12").Value = 0.1798561151079137
    $ws.Range("K12").Value = 0.01438848920863309
    $ws.Range("L12").Value = 0.04316546762589928
    $ws.Range("S12").Value = 0.02877697841726619
    $ws.Range("G13").Value = 0.66
    $ws.Range("J13").Value = 0.22
    $ws.Range("S13").Value = 0.12
    $ws.Range("J14").Value = 1
    $ws.Range("F15").Value = 0.02051282051282051
    $ws.Range("H15").Value = 0.1846153846153846
    $ws.Range("I15").Value = 0.05128205128205128
    $ws.Range("J15").Value = 0.358974358974359
    $ws.Range("K15").Value = 0.06666666666666667
    $ws.Range("M15").Value = 0.01538461538461539
    $ws.Range("O15").Value = 0.07692307692307693
    $ws.Range("S15").Value = 0.2256410256410256
    $ws.Range("F16").Value = 0.0308641975308642
    $ws.Range("H16").Value = 0.1975308641975309
    $ws.Range("I16").Value = 0.1049382716049383
    $ws.Range("J16").Value = 0.3888888888888889
    $ws.Range("K16").Value = 0.1172839506172839
    $ws.Range("M16").Value = 0.02469135802469136
    $ws.Range("O16").Value = 0.02469135802469136
    $ws.Range("S16").Value = 0.1111111111111111
    $ws.Range("F17").Value = 0.004454342984409799
    $ws.Range("H17").Value = 0.155902004454343
    $ws.Range("I17").Value = 0.09131403118040089
    $ws.Range("J17").Value = 0.4543429844097995
    $ws.Range("K17").Value = 0.08685968819599109
    $ws.Range("M17").Value = 0.0200445434298441
    $ws.Range("O17").Value = 0.0801781737193764
    $ws.Range("S17").Value = 0.1069042316258352
    $ws.Range("F18").Value = 0.02985074626865672
    $ws.Range("H18").Value = 0.1890547263681592
    $ws.Range("I18").Value = 0.07960199004975124
    $ws.Range("J18").Value = 0.4029850746268657
    $ws.Range("K18").Value = 0.0845771144278607
    $ws.Range("M18").Value = 0.01492537313432836
    $ws.Range("O18").Value = 0.0945273631840796
    $ws.Range("S18").Value = 0.1044776119402985
    $ws.Range("F19").Value = 0.01546860782529572
    $ws.Range("H19").Value = 0.2202001819836215
    $ws.Range("I19").Value = 0.09463148316651501
    $ws.Range("J19").Value = 0.362147406733394
    $ws.Range("K19").Value = 0.08371246587807098
    $ws.Range("M19").Value = 0.0272975432211101
    $ws.Range("N19").Value = 0.0009099181073703367
    $ws.Range("O19").Value = 0.07097361237488627
    $ws.Range("S19").Value = 0.1246587807097361
